# Update Name of Algo
# Apply updated imputed values produced by the KNN algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 6.449
$ws.Range("E4").Value  = 12.718
$ws.Range("E5").Value  = 13.499
$ws.Range("B6").Value  = 7.295999999999999
$ws.Range("B7").Value  = 6.612
$ws.Range("E8").Value  = 13.718
$ws.Range("B16").Value = 6.782999999999999
$ws.Range("E16").Value = 13.159
$ws.Range("B20").Value = 6.145
$ws.Range("E22").Value = 13.495
